$wb = $excel.ActiveWorkbook

# Sheet "deploy_amounts": update row 2 allocation amounts
$wsAmounts = $wb.Worksheets.Item("deploy_amounts")
$wsAmounts.Range("A2").Value = 80000
$wsAmounts.Range("T2").Value = 0
$wsAmounts.Range("W2").Value = 80000
$wsAmounts.Range("AQ2").Value = 0
$wsAmounts.Range("BC2").Value = 80000

# Sheet "deploy_bins": update which rows are flagged (binary indicator)
$wsBins = $wb.Worksheets.Item("deploy_bins")
$wsBins.Range("A21").Value = 0
$wsBins.Range("A24").Value = 1
$wsBins.Range("A44").Value = 0
$wsBins.Range("A56").Value = 1
